# Auto-generated Excel COM-interop script
# Applies the 'Updated cryptos list' data refresh to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.652.91'
$ws.Range("E2").Value = '  +3.48%  '
$ws.Range("D3").Value = '2.406.73'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.30%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.571'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.26%  '
$ws.Range("E9").Value = '  +6.56%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.83'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.364'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.150'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.70'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.91%  '
$ws.Range("D14").Value = '2.834.05'
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("D15").Value = '59.547.97'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000139'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.94%  '
$ws.Range("D17").Value = '2.409.05'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.43%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '336.69'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.59%  '
$ws.Range("E21").Value = '  +5.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.91'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.28%  '
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.37'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.43%  '
$ws.Range("E28").Value = '  +2.93%  '
$ws.Range("D29").Value = '0.0₃0763'
$ws.Range("E29").Value = '  +4.88%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.22'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.78'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.31'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '40.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.68%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.422'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.18%  '
$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '300.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.08%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.42'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0962'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0525'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.572'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.51%  '
$ws.Range("E48").Value = '  +4.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.399'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.86%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("E51").Value = '  +4.90%  '
